$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting C:T to D:U
$ws.Columns("C").Insert()

# Set header for the newly inserted column C1
$ws.Range("C1").Value = "Unnamed: 0"

# Rename existing B1 header from "Unnamed: 0" to "Unnamed: 0.1"
$ws.Range("B1").Value = "Unnamed: 0.1"
